# Update "想去人数" (F column) counts on both the "展览" sheet and the
# aggregated "全部类型" sheet, per the latest gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(4, 6).Value = 1586
$ws1.Cells.Item(5, 6).Value = 174
$ws1.Cells.Item(8, 6).Value = 197
$ws1.Cells.Item(9, 6).Value = 774
$ws1.Cells.Item(10, 6).Value = 1061
$ws1.Cells.Item(12, 6).Value = 367
$ws1.Cells.Item(14, 6).Value = 524
$ws1.Cells.Item(15, 6).Value = 32
$ws1.Cells.Item(16, 6).Value = 6584
$ws1.Cells.Item(17, 6).Value = 31
$ws1.Cells.Item(20, 6).Value = 171
$ws1.Cells.Item(22, 6).Value = 1022
$ws1.Cells.Item(23, 6).Value = 15743
$ws1.Cells.Item(24, 6).Value = 1556
$ws1.Cells.Item(25, 6).Value = 22
$ws1.Cells.Item(26, 6).Value = 305
$ws1.Cells.Item(27, 6).Value = 157
$ws1.Cells.Item(29, 6).Value = 11163
$ws1.Cells.Item(30, 6).Value = 799
$ws1.Cells.Item(31, 6).Value = 4380
$ws1.Cells.Item(32, 6).Value = 272
$ws1.Cells.Item(35, 6).Value = 311
$ws1.Cells.Item(36, 6).Value = 132

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(4, 6).Value = 1586
$ws4.Cells.Item(5, 6).Value = 174
$ws4.Cells.Item(9, 6).Value = 197
$ws4.Cells.Item(10, 6).Value = 775
$ws4.Cells.Item(12, 6).Value = 1061
$ws4.Cells.Item(14, 6).Value = 367
$ws4.Cells.Item(16, 6).Value = 524
$ws4.Cells.Item(18, 6).Value = 32
$ws4.Cells.Item(19, 6).Value = 6584
$ws4.Cells.Item(20, 6).Value = 31
$ws4.Cells.Item(23, 6).Value = 171
$ws4.Cells.Item(26, 6).Value = 1022
$ws4.Cells.Item(27, 6).Value = 15743
$ws4.Cells.Item(28, 6).Value = 1556
$ws4.Cells.Item(29, 6).Value = 22
$ws4.Cells.Item(30, 6).Value = 305
$ws4.Cells.Item(31, 6).Value = 157
$ws4.Cells.Item(34, 6).Value = 11163
$ws4.Cells.Item(35, 6).Value = 799
$ws4.Cells.Item(36, 6).Value = 4380
$ws4.Cells.Item(37, 6).Value = 272
$ws4.Cells.Item(40, 6).Value = 311
$ws4.Cells.Item(41, 6).Value = 132
